# Auto-generated edit script: updates FFXIV leve-profit market-price
# snapshot values (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, matching the 'chore: update Sheets via scheduled runner' commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 3771.353
$ws.Cells.Item(64, 9).Value = 4192.8
$ws.Cells.Item(64, 10).Value = 3595.75
$ws.Cells.Item(64, 11).Value = 4192.8
$ws.Cells.Item(64, 12).Value = 3595.75
$ws.Cells.Item(64, 13).Value = -3944.8
$ws.Cells.Item(64, 14).Value = -4091.75

# Row 67
$ws.Cells.Item(67, 8).Value = 3771.353
$ws.Cells.Item(67, 9).Value = 4192.8
$ws.Cells.Item(67, 10).Value = 3595.75
$ws.Cells.Item(67, 11).Value = 4192.8
$ws.Cells.Item(67, 12).Value = 3595.75
$ws.Cells.Item(67, 13).Value = -3334.8
$ws.Cells.Item(67, 14).Value = -5311.75

# Row 100
$ws.Cells.Item(100, 8).Value = 575.27905
$ws.Cells.Item(100, 9).Value = 492.76315
$ws.Cells.Item(100, 10).Value = 1202.4
$ws.Cells.Item(100, 11).Value = 492.76315
$ws.Cells.Item(100, 12).Value = 1202.4
$ws.Cells.Item(100, 13).Value = 48.23685
$ws.Cells.Item(100, 14).Value = -2284.4

# Row 106
$ws.Cells.Item(106, 8).Value = 3302.2666
$ws.Cells.Item(106, 9).Value = 3335.6667
$ws.Cells.Item(106, 10).Value = 3168.6667
$ws.Cells.Item(106, 11).Value = 3335.6667
$ws.Cells.Item(106, 12).Value = 3168.6667
$ws.Cells.Item(106, 13).Value = -2704.6667
$ws.Cells.Item(106, 14).Value = -4430.6667

# Row 111
$ws.Cells.Item(111, 8).Value = 3135.4285
$ws.Cells.Item(111, 9).Value = 4037
$ws.Cells.Item(111, 10).Value = 1933.3334
$ws.Cells.Item(111, 11).Value = 12111
$ws.Cells.Item(111, 12).Value = 5800.0002
$ws.Cells.Item(111, 13).Value = -9044
$ws.Cells.Item(111, 14).Value = -11934.0002

# Row 125
$ws.Cells.Item(125, 8).Value = 1620.5555
$ws.Cells.Item(125, 9).Value = 1949.8572
$ws.Cells.Item(125, 10).Value = 468
$ws.Cells.Item(125, 11).Value = 17548.7148
$ws.Cells.Item(125, 12).Value = 4212
$ws.Cells.Item(125, 13).Value = -15088.7148
$ws.Cells.Item(125, 14).Value = -9132

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 10092.454
$ws.Cells.Item(2, 9).Value = 1061.6
$ws.Cells.Item(2, 11).Value = 1061.6
$ws.Cells.Item(2, 13).Value = -948.5999999999999

# Row 28
$ws.Cells.Item(28, 8).Value = 5027.143
$ws.Cells.Item(28, 9).Value = 5027.143
$ws.Cells.Item(28, 11).Value = 5027.143
$ws.Cells.Item(28, 13).Value = -4835.143

# Row 31
$ws.Cells.Item(31, 8).Value = 17166
$ws.Cells.Item(31, 9).Value = 8124
$ws.Cells.Item(31, 10).Value = 35250
$ws.Cells.Item(31, 11).Value = 8124
$ws.Cells.Item(31, 12).Value = 35250
$ws.Cells.Item(31, 13).Value = -7830
$ws.Cells.Item(31, 14).Value = -35838

# Row 32
$ws.Cells.Item(32, 8).Value = 4170.596
$ws.Cells.Item(32, 9).Value = 3970.6445
$ws.Cells.Item(32, 11).Value = 3970.6445
$ws.Cells.Item(32, 13).Value = -3683.6445

# Row 61
$ws.Cells.Item(61, 8).Value = 2027.75
$ws.Cells.Item(61, 9).Value = 1560
$ws.Cells.Item(61, 10).Value = 2682.6
$ws.Cells.Item(61, 11).Value = 1560
$ws.Cells.Item(61, 12).Value = 2682.6
$ws.Cells.Item(61, 13).Value = -1348
$ws.Cells.Item(61, 14).Value = -3106.6

# Row 99
$ws.Cells.Item(99, 8).Value = 5027.143
$ws.Cells.Item(99, 9).Value = 5027.143
$ws.Cells.Item(99, 11).Value = 5027.143
$ws.Cells.Item(99, 13).Value = -2032.143

# Row 102
$ws.Cells.Item(102, 8).Value = 20835140
$ws.Cells.Item(102, 9).Value = 23810874
$ws.Cells.Item(102, 11).Value = 23810874
$ws.Cells.Item(102, 13).Value = -23809252

# Row 110
$ws.Cells.Item(110, 8).Value = 1726.3
$ws.Cells.Item(110, 9).Value = 825
$ws.Cells.Item(110, 11).Value = 825
$ws.Cells.Item(110, 13).Value = 1220

# Row 116
$ws.Cells.Item(116, 8).Value = 10092.454
$ws.Cells.Item(116, 9).Value = 1061.6
$ws.Cells.Item(116, 11).Value = 1061.6
$ws.Cells.Item(116, 13).Value = 1232.4

# Row 136
$ws.Cells.Item(136, 8).Value = 2027.75
$ws.Cells.Item(136, 9).Value = 1560
$ws.Cells.Item(136, 10).Value = 2682.6
$ws.Cells.Item(136, 11).Value = 4680
$ws.Cells.Item(136, 12).Value = 8047.799999999999
$ws.Cells.Item(136, 13).Value = -2130
$ws.Cells.Item(136, 14).Value = -13147.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 10092.454
$ws.Cells.Item(3, 9).Value = 1061.6
$ws.Cells.Item(3, 11).Value = 1061.6
$ws.Cells.Item(3, 13).Value = -947.5999999999999

# Row 20
$ws.Cells.Item(20, 8).Value = 2380.5625
$ws.Cells.Item(20, 9).Value = 2255.5557
$ws.Cells.Item(20, 11).Value = 2255.5557
$ws.Cells.Item(20, 13).Value = -2008.5557

# Row 86
$ws.Cells.Item(86, 8).Value = 2637.9546
$ws.Cells.Item(86, 9).Value = 2782.3635
$ws.Cells.Item(86, 10).Value = 2204.7273
$ws.Cells.Item(86, 11).Value = 2782.3635
$ws.Cells.Item(86, 12).Value = 2204.7273
$ws.Cells.Item(86, 13).Value = -1659.3635
$ws.Cells.Item(86, 14).Value = -4450.7273

# Row 89
$ws.Cells.Item(89, 8).Value = 2637.9546
$ws.Cells.Item(89, 9).Value = 2782.3635
$ws.Cells.Item(89, 10).Value = 2204.7273
$ws.Cells.Item(89, 11).Value = 13911.8175
$ws.Cells.Item(89, 12).Value = 11023.6365
$ws.Cells.Item(89, 13).Value = -8295.817499999999
$ws.Cells.Item(89, 14).Value = -22255.6365

# Row 99
$ws.Cells.Item(99, 8).Value = 29413142
$ws.Cells.Item(99, 9).Value = 38462860
$ws.Cells.Item(99, 10).Value = 1552.5
$ws.Cells.Item(99, 11).Value = 38462860
$ws.Cells.Item(99, 12).Value = 1552.5
$ws.Cells.Item(99, 13).Value = -38461362
$ws.Cells.Item(99, 14).Value = -4548.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 100001140
$ws.Cells.Item(16, 9).Value = 111112210
$ws.Cells.Item(16, 10).Value = 1500
$ws.Cells.Item(16, 11).Value = 111112210
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = -111111923
$ws.Cells.Item(16, 14).Value = -2074

# Row 50
$ws.Cells.Item(50, 8).Value = 18541.5
$ws.Cells.Item(50, 9).Value = 10083
$ws.Cells.Item(50, 10).Value = 27000
$ws.Cells.Item(50, 11).Value = 10083
$ws.Cells.Item(50, 12).Value = 27000
$ws.Cells.Item(50, 13).Value = -9458
$ws.Cells.Item(50, 14).Value = -28250

# Row 51
$ws.Cells.Item(51, 8).Value = 16697
$ws.Cells.Item(51, 9).Value = 1091
$ws.Cells.Item(51, 10).Value = 24500
$ws.Cells.Item(51, 11).Value = 1091
$ws.Cells.Item(51, 12).Value = 24500
$ws.Cells.Item(51, 13).Value = -355
$ws.Cells.Item(51, 14).Value = -25972

# Row 60
$ws.Cells.Item(60, 8).Value = 9076.695
$ws.Cells.Item(60, 10).Value = 12096.5625
$ws.Cells.Item(60, 12).Value = 12096.5625
$ws.Cells.Item(60, 14).Value = -13118.5625

# Row 61
$ws.Cells.Item(61, 8).Value = 16697
$ws.Cells.Item(61, 9).Value = 1091
$ws.Cells.Item(61, 10).Value = 24500
$ws.Cells.Item(61, 11).Value = 1091
$ws.Cells.Item(61, 12).Value = 24500
$ws.Cells.Item(61, 13).Value = -743
$ws.Cells.Item(61, 14).Value = -25196

# Row 62
$ws.Cells.Item(62, 8).Value = 100012500
$ws.Cells.Item(62, 9).Value = 25000
$ws.Cells.Item(62, 11).Value = 25000
$ws.Cells.Item(62, 13).Value = -24376

# Row 65
$ws.Cells.Item(65, 8).Value = 100012500
$ws.Cells.Item(65, 9).Value = 25000
$ws.Cells.Item(65, 11).Value = 125000
$ws.Cells.Item(65, 13).Value = -121880

# Row 99
$ws.Cells.Item(99, 8).Value = 1865.2667
$ws.Cells.Item(99, 9).Value = 1807
$ws.Cells.Item(99, 10).Value = 1981.8
$ws.Cells.Item(99, 11).Value = 1807
$ws.Cells.Item(99, 12).Value = 1981.8
$ws.Cells.Item(99, 13).Value = -309
$ws.Cells.Item(99, 14).Value = -4977.8

# Row 113
$ws.Cells.Item(113, 8).Value = 100001140
$ws.Cells.Item(113, 9).Value = 111112210
$ws.Cells.Item(113, 10).Value = 1500
$ws.Cells.Item(113, 11).Value = 111112210
$ws.Cells.Item(113, 12).Value = 1500
$ws.Cells.Item(113, 13).Value = -111110040
$ws.Cells.Item(113, 14).Value = -5840

# Row 122
$ws.Cells.Item(122, 8).Value = 1043.5555
$ws.Cells.Item(122, 9).Value = 999
$ws.Cells.Item(122, 10).Value = 1400
$ws.Cells.Item(122, 11).Value = 2997
$ws.Cells.Item(122, 12).Value = 4200
$ws.Cells.Item(122, 13).Value = -547
$ws.Cells.Item(122, 14).Value = -9100

# Row 126
$ws.Cells.Item(126, 8).Value = 1865.2667
$ws.Cells.Item(126, 9).Value = 1807
$ws.Cells.Item(126, 10).Value = 1981.8
$ws.Cells.Item(126, 11).Value = 5421
$ws.Cells.Item(126, 12).Value = 5945.4
$ws.Cells.Item(126, 13).Value = -2951
$ws.Cells.Item(126, 14).Value = -10885.4

# Row 134
$ws.Cells.Item(134, 8).Value = 23811412
$ws.Cells.Item(134, 9).Value = 30304980
$ws.Cells.Item(134, 10).Value = 1662.3334
$ws.Cells.Item(134, 11).Value = 90914940
$ws.Cells.Item(134, 12).Value = 4987.0002
$ws.Cells.Item(134, 13).Value = -90912405
$ws.Cells.Item(134, 14).Value = -10057.0002

# Row 141
$ws.Cells.Item(141, 8).Value = 34549.09
$ws.Cells.Item(141, 10).Value = 34549.09
$ws.Cells.Item(141, 12).Value = 34549.09
$ws.Cells.Item(141, 14).Value = -44909.09

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 15385866
$ws.Cells.Item(131, 10).Value = 1401.5536
$ws.Cells.Item(131, 12).Value = 4204.6608
$ws.Cells.Item(131, 14).Value = -14284.6608

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 779.3570999999999
$ws.Cells.Item(107, 9).Value = 930.25
$ws.Cells.Item(107, 10).Value = 578.1667
$ws.Cells.Item(107, 11).Value = 930.25
$ws.Cells.Item(107, 12).Value = 578.1667
$ws.Cells.Item(107, 13).Value = 989.75
$ws.Cells.Item(107, 14).Value = -4418.1667

# Row 122
$ws.Cells.Item(122, 8).Value = 168688.67
$ws.Cells.Item(122, 9).Value = 2385.4285
$ws.Cells.Item(122, 11).Value = 7156.2855
$ws.Cells.Item(122, 13).Value = -4706.2855

# Row 126
$ws.Cells.Item(126, 8).Value = 2418.0833
$ws.Cells.Item(126, 9).Value = 2224.111
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 6672.333
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -4202.333
$ws.Cells.Item(126, 14).Value = -13940

# Row 132
$ws.Cells.Item(132, 8).Value = 3297.8125
$ws.Cells.Item(132, 9).Value = 3386.6667
$ws.Cells.Item(132, 10).Value = 3244.5
$ws.Cells.Item(132, 11).Value = 10160.0001
$ws.Cells.Item(132, 12).Value = 9733.5
$ws.Cells.Item(132, 13).Value = -7630.000100000001
$ws.Cells.Item(132, 14).Value = -14793.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 1627.619
$ws.Cells.Item(22, 9).Value = 1598.8125
$ws.Cells.Item(22, 10).Value = 1719.8
$ws.Cells.Item(22, 11).Value = 1598.8125
$ws.Cells.Item(22, 12).Value = 1719.8
$ws.Cells.Item(22, 13).Value = -1303.8125
$ws.Cells.Item(22, 14).Value = -2309.8

# Row 27
$ws.Cells.Item(27, 8).Value = 1627.619
$ws.Cells.Item(27, 9).Value = 1598.8125
$ws.Cells.Item(27, 10).Value = 1719.8
$ws.Cells.Item(27, 11).Value = 1598.8125
$ws.Cells.Item(27, 12).Value = 1719.8
$ws.Cells.Item(27, 13).Value = -1491.8125
$ws.Cells.Item(27, 14).Value = -1933.8

# Row 46
$ws.Cells.Item(46, 8).Value = 5997
$ws.Cells.Item(46, 9).Value = 1992.5
$ws.Cells.Item(46, 11).Value = 1992.5
$ws.Cells.Item(46, 13).Value = -1804.5

# Row 61
$ws.Cells.Item(61, 8).Value = 916.6667
$ws.Cells.Item(61, 9).Value = 700
$ws.Cells.Item(61, 11).Value = 700
$ws.Cells.Item(61, 13).Value = -498

# Row 68
$ws.Cells.Item(68, 8).Value = 1809.7778
$ws.Cells.Item(68, 9).Value = 1570
$ws.Cells.Item(68, 11).Value = 1570
$ws.Cells.Item(68, 13).Value = -821

# Row 71
$ws.Cells.Item(71, 8).Value = 1809.7778
$ws.Cells.Item(71, 9).Value = 1570
$ws.Cells.Item(71, 11).Value = 7850
$ws.Cells.Item(71, 13).Value = -4106

# Row 113
$ws.Cells.Item(113, 8).Value = 916.6667
$ws.Cells.Item(113, 9).Value = 700
$ws.Cells.Item(113, 11).Value = 700
$ws.Cells.Item(113, 13).Value = 1470

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Cells.Item(63, 8).Value = 21750
$ws.Cells.Item(63, 10).Value = 21750
$ws.Cells.Item(63, 12).Value = 21750
$ws.Cells.Item(63, 14).Value = -22998

# Row 66
$ws.Cells.Item(66, 8).Value = 21750
$ws.Cells.Item(66, 10).Value = 21750
$ws.Cells.Item(66, 12).Value = 65250
$ws.Cells.Item(66, 14).Value = -71490

# Row 122
$ws.Cells.Item(122, 8).Value = 10402027
$ws.Cells.Item(122, 9).Value = 10402027
$ws.Cells.Item(122, 10).Value = 10402027
$ws.Cells.Item(122, 11).Value = 31206081
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -31206081
$ws.Cells.Item(122, 14).ClearContents()
